$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0).Date

# Data rows run from row 2 through row 143 (column A holds the "Beteckning" id
# for each row); update the "Förändrad" date in column C for each of them.
$lastRow = 143
for ($r = 2; $r -le $lastRow; $r++) {
    $idCell = $ws.Cells.Item($r, 1)
    if ($idCell.Value -ne $null -and $idCell.Value -ne "") {
        $ws.Cells.Item($r, 3).Value = $newDate
    }
}
